$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. '27.035.08', '311.06') that
# must stay plain text, matching the source inline-string cells. Setting
# .Value directly on such strings makes Excel's COM layer coerce them to
# real numbers, so we force text mode first and restore the default
# 'Normal' style afterwards (keeps t="s"/text without leaving a stray
# number-format style applied to the cell).
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "27.035.08"
Set-TextValue "D3" "1.821.04"
$ws.Range("E4").Value = "  -0.51%  "
Set-TextValue "D5" "311.06"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("E6").Value = "  -0.46%  "
Set-TextValue "D7" "0.4492"
$ws.Range("E7").Value = "  +5.14%  "
Set-TextValue "D8" "0.3691"
$ws.Range("E8").Value = "  +0.11%  "
Set-TextValue "D9" "0.07289"
$ws.Range("E9").Value = "  +0.62%  "
Set-TextValue "D10" "0.8571"
$ws.Range("E10").Value = "  -0.61%  "
Set-TextValue "D11" "20.77"
$ws.Range("E11").Value = "  -1.13%  "
Set-TextValue "D12" "1.823.87"
$ws.Range("E12").Value = "  -0.07%  "
Set-TextValue "D13" "6.642"
$ws.Range("E13").Value = "  -1.34%  "
Set-TextValue "D14" "92.48"
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D15" "5.337"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D16" "0.07097"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("E17").Value = "  -0.46%  "
Set-TextValue "D18" "0.000008790"
$ws.Range("E19").Value = "  -0.47%  "
Set-TextValue "D20" "14.98"
$ws.Range("E20").Value = "  -0.67%  "
Set-TextValue "D21" "27.049.49"
$ws.Range("E21").Value = "  -0.74%  "
Set-TextValue "D22" "5.163"
$ws.Range("E22").Value = "  +0.40%  "
Set-TextValue "D23" "10.93"
$ws.Range("E23").Value = "  +0.44%  "
Set-TextValue "D24" "1.993"
$ws.Range("E24").Value = "  -0.59%  "
Set-TextValue "D25" "151.68"
$ws.Range("E25").Value = "  -0.95%  "
Set-TextValue "D26" "2.236"
$ws.Range("E26").Value = "  +4.89%  "
Set-TextValue "D27" "18.47"
$ws.Range("E27").Value = "  +0.58%  "
Set-TextValue "D28" "5.249"
$ws.Range("E28").Value = "  +0.29%  "
Set-TextValue "D29" "116.50"
$ws.Range("E29").Value = "  +0.16%  "
Set-TextValue "D30" "0.08866"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D31" "0.7550"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D32" "1.183"
$ws.Range("E32").Value = "  -1.17%  "
Set-TextValue "D33" "2.962"
$ws.Range("E33").Value = "  +4.38%  "
Set-TextValue "D34" "4.454"
$ws.Range("E34").Value = "  +0.41%  "
Set-TextValue "D35" "1.001"
$ws.Range("E35").Value = "  -0.53%  "
Set-TextValue "D36" "1.093"
$ws.Range("E36").Value = "  -1.84%  "
Set-TextValue "D37" "0.01969"
$ws.Range("E37").Value = "  -0.11%  "
Set-TextValue "D38" "0.05236"
$ws.Range("E38").Value = "  -0.67%  "
Set-TextValue "D39" "0.5329"
$ws.Range("E39").Value = "  +5.69%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D40" "2.888"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "7.148"
$ws.Range("E41").Value = "  +0.18%  "
Set-TextValue "D42" "0.1704"
$ws.Range("E42").Value = "  +0.39%  "
Set-TextValue "D43" "0.5263"
$ws.Range("E43").Value = "  +11.03%  "
Set-TextValue "D44" "8.524"
$ws.Range("E44").Value = "  -1.42%  "
Set-TextValue "D45" "10.58"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  +8.67%  "
Set-TextValue "D47" "105.61"
$ws.Range("E47").Value = "  -1.73%  "
Set-TextValue "D48" "1.001"
$ws.Range("E48").Value = "  -0.48%  "
Set-TextValue "D49" "1.667"
$ws.Range("E49").Value = "  +0.43%  "
Set-TextValue "D50" "0.06385"
$ws.Range("E50").Value = "  +0.23%  "
Set-TextValue "D51" "0.9196"
$ws.Range("E51").Value = "  +0.38%  "
